$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column was populated with a mangled value derived from the
# source file name (e.g. "6-13-2007-08") instead of the actual game date
# ("2008-06-13"). NBA "season" stats files are named after the season
# (2007-08) plus the in-season file date, which is one day off from the
# real game date once normalized to ISO form - this restores the correct
# value for every data row.

$oldValue = "6-13-2007-08"
$newValue = "2008-06-13"

# Locate the "Date" header cell (column BF in this workbook) dynamically
# rather than hard-coding the column letter.
$headerRow = 1
$dateHeaderCell = $ws.Rows.Item($headerRow).Find("Date")

$used = $ws.UsedRange
$lastRow = $used.Row() + $used.Rows.Count() - 1
$dateCol = $dateHeaderCell.Column()

for ($row = $headerRow + 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    $val = $cell.Value()
    if ($val -eq $oldValue) {
        # Force text formatting before assigning, otherwise Excel parses
        # the ISO-looking string and silently turns it into a date serial
        # number instead of keeping it as literal text.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        # Restore the default style so only the value itself changes.
        $cell.Style = "Normal"
    }
}
